{"js": "// Update the date heading (first paragraph) and the answers in the\n// 20x5 addition/subtraction practice table.\nconst body = context.document.body;\n\n// --- 1. Update the date/weekday heading (first paragraph in the body) ---\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst titleParagraph = paragraphs.items[0];\ntitleParagraph.insertText(\"2024-12-12 Thursday\", Word.InsertLocation.replace);\n\n// --- 2. Update every answer cell in the practice table ---\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount, values\");\nawait context.sync();\n\n// New answers, in row-major order (row 1 col 1..5, row 2 col 1..5, ...),\n// matching the table's existing 20-row x 5-column layout.\nconst newValues = [\n  [\n    \"32+36=68\",\n    \"50-21=29\",\n    \"13+37=50\",\n    \"60+19=79\",\n    \"98-70=28\"\n  ],\n  [\n    \"97-50=47\",\n    \"86-74=12\",\n    \"11+26=37\",\n    \"5+21=26\",\n    \"67+24=91\"\n  ],\n  [\n    \"86-73=13\",\n    \"87-76=11\",\n    \"99-23=76\",\n    \"37+15=52\",\n    \"45-39=6\"\n  ],\n  [\n    \"2-1=1\",\n    \"69-59=10\",\n    \"80-13=67\",\n    \"23+36=59\",\n    \"82-17=65\"\n  ],\n  [\n    \"4+38=42\",\n    \"43+12=55\",\n    \"53-33=20\",\n    \"4+4=8\",\n    \"19-3=16\"\n  ],\n  [\n    \"41+5=46\",\n    \"94-89=5\",\n    \"92-32=60\",\n    \"88-19=69\",\n    \"73-43=30\"\n  ],\n  [\n    \"96-31=65\",\n    \"85-71=14\",\n    \"90-44=46\",\n    \"80-1=79\",\n    \"7+32=39\"\n  ],\n  [\n    \"57+41=98\",\n    \"7+86=93\",\n    \"21+65=86\",\n    \"67-47=20\",\n    \"30+63=93\"\n  ],\n  [\n    \"79-39=40\",\n    \"50+48=98\",\n    \"56-20=36\",\n    \"66+14=80\",\n    \"50+11=61\"\n  ],\n  [\n    \"58-31=27\",\n    \"15+54=69\",\n    \"96-18=78\",\n    \"73-56=17\",\n    \"36+16=52\"\n  ],\n  [\n    \"80-74=6\",\n    \"84+14=98\",\n    \"15+33=48\",\n    \"85-82=3\",\n    \"84-48=36\"\n  ],\n  [\n    \"31+19=50\",\n    \"18+58=76\",\n    \"88+8=96\",\n    \"7-1=6\",\n    \"5+10=15\"\n  ],\n  [\n    \"20+25=45\",\n    \"35+50=85\",\n    \"50-27=23\",\n    \"30+39=69\",\n    \"51-15=36\"\n  ],\n  [\n    \"19+43=62\",\n    \"50+47=97\",\n    \"98-12=86\",\n    \"8+59=67\",\n    \"67-16=51\"\n  ],\n  [\n    \"41+37=78\",\n    \"38-25=13\",\n    \"65+31=96\",\n    \"17+1=18\",\n    \"8+1=9\"\n  ],\n  [\n    \"35-14=21\",\n    \"23+39=62\",\n    \"44-27=17\",\n    \"17+42=59\",\n    \"86-81=5\"\n  ],\n  [\n    \"14+30=44\",\n    \"37+33=70\",\n    \"34-18=16\",\n    \"57+16=73\",\n    \"24+15=39\"\n  ],\n  [\n    \"34-23=11\",\n    \"78-33=45\",\n    \"43+19=62\",\n    \"8+10=18\",\n    \"39-27=12\"\n  ],\n  [\n    \"20-16=4\",\n    \"81-22=59\",\n    \"10+20=30\",\n    \"98-18=80\",\n    \"90-47=43\"\n  ],\n  [\n    \"6+54=60\",\n    \"92-82=10\",\n    \"83-8=75\",\n    \"22+8=30\",\n    \"21-11=10\"\n  ]\n];\n\nif (table.rowCount === newValues.length &&\n    table.values.length > 0 &&\n    table.values[0].length === newValues[0].length) {\n  // Bulk-write the whole table in one shot.\n  table.values = newValues;\n} else {\n  // Fallback: write cell-by-cell if the table shape ever changes.\n  for (let r = 0; r < newValues.length; r++) {\n    for (let c = 0; c < newValues[r].length; c++) {\n      const cell = table.getCell(r, c);\n      cell.body.insertText(newValues[r][c], Word.InsertLocation.replace);\n    }\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date heading (first paragraph) and the 20x5 practice-problem table.\n$d = $word.ActiveDocument\n\n# New date string for the title paragraph.\n$d.Paragraphs.Item(1).Range.Text = \"2024-12-12 Thursday\"\n\n# New answers for each of the 100 table cells, in row-major order\n# (row 1 col 1..5, row 2 col 1..5, ...).\n$values = @(\n    \"32+36=68\",\n    \"50-21=29\",\n    \"13+37=50\",\n    \"60+19=79\",\n    \"98-70=28\",\n    \"97-50=47\",\n    \"86-74=12\",\n    \"11+26=37\",\n    \"5+21=26\",\n    \"67+24=91\",\n    \"86-73=13\",\n    \"87-76=11\",\n    \"99-23=76\",\n    \"37+15=52\",\n    \"45-39=6\",\n    \"2-1=1\",\n    \"69-59=10\",\n    \"80-13=67\",\n    \"23+36=59\",\n    \"82-17=65\",\n    \"4+38=42\",\n    \"43+12=55\",\n    \"53-33=20\",\n    \"4+4=8\",\n    \"19-3=16\",\n    \"41+5=46\",\n    \"94-89=5\",\n    \"92-32=60\",\n    \"88-19=69\",\n    \"73-43=30\",\n    \"96-31=65\",\n    \"85-71=14\",\n    \"90-44=46\",\n    \"80-1=79\",\n    \"7+32=39\",\n    \"57+41=98\",\n    \"7+86=93\",\n    \"21+65=86\",\n    \"67-47=20\",\n    \"30+63=93\",\n    \"79-39=40\",\n    \"50+48=98\",\n    \"56-20=36\",\n    \"66+14=80\",\n    \"50+11=61\",\n    \"58-31=27\",\n    \"15+54=69\",\n    \"96-18=78\",\n    \"73-56=17\",\n    \"36+16=52\",\n    \"80-74=6\",\n    \"84+14=98\",\n    \"15+33=48\",\n    \"85-82=3\",\n    \"84-48=36\",\n    \"31+19=50\",\n    \"18+58=76\",\n    \"88+8=96\",\n    \"7-1=6\",\n    \"5+10=15\",\n    \"20+25=45\",\n    \"35+50=85\",\n    \"50-27=23\",\n    \"30+39=69\",\n    \"51-15=36\",\n    \"19+43=62\",\n    \"50+47=97\",\n    \"98-12=86\",\n    \"8+59=67\",\n    \"67-16=51\",\n    \"41+37=78\",\n    \"38-25=13\",\n    \"65+31=96\",\n    \"17+1=18\",\n    \"8+1=9\",\n    \"35-14=21\",\n    \"23+39=62\",\n    \"44-27=17\",\n    \"17+42=59\",\n    \"86-81=5\",\n    \"14+30=44\",\n    \"37+33=70\",\n    \"34-18=16\",\n    \"57+16=73\",\n    \"24+15=39\",\n    \"34-23=11\",\n    \"78-33=45\",\n    \"43+19=62\",\n    \"8+10=18\",\n    \"39-27=12\",\n    \"20-16=4\",\n    \"81-22=59\",\n    \"10+20=30\",\n    \"98-18=80\",\n    \"90-47=43\",\n    \"6+54=60\",\n    \"92-82=10\",\n    \"83-8=75\",\n    \"22+8=30\",\n    \"21-11=10\"\n)\n\n$t = $d.Tables.Item(1)\n$rows = $t.Rows.Count\n$cols = $t.Columns.Count\n\nfor ($i = 0; $i -lt $values.Length; $i++) {\n    $row = [int][math]::Floor($i / $cols) + 1\n    $col = ($i % $cols) + 1\n    $cell = $t.Cell($row, $col)\n    $cell.Range.Text = $values[$i]\n}\n"}
